$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 0.04462476567923779
$ws.Range("I2").Value = 0.04462476567923779
$ws.Range("J2").Value = 0.04064155697093108
$ws.Range("K2").Value = 0.04064155697093108
$ws.Range("L2").Value = 7.466497182163403
$ws.Range("M2").Value = "[0.32306204900709545, 14.60993231531971]"
$ws.Range("N2").Value = 0.04088902655447568
$ws.Range("O2").Value = 0.04088902655447568
$ws.Range("P2").Value = -1.371105502467618
$ws.Range("Q2").Value = "[-2.7422110049352346, -8.881784197001252e-16]"
$ws.Range("R2").Value = 0.0500000004411274
$ws.Range("S2").Value = 0.0500000004411274
$ws.Range("T2").Value = 11.25306146723949
$ws.Range("U2").Value = "[7.230002838331853, 15.276120096147135]"
$ws.Range("V2").Value = [double]"1.08840863832782e-06"
$ws.Range("W2").Value = [double]"1.08840863832782e-06"
$ws.Range("X2").Value = 4.800800800800801
$ws.Range("Y2").Value = [double]"1.77635683940025e-15"
$ws.Range("Z2").Value = 9.601601601601599
$ws.Range("H3").Value = 0.3172302047542273
$ws.Range("I3").Value = 0.3172302047542273
$ws.Range("J3").Value = 0.008403349655002557
$ws.Range("K3").Value = 0.008403349655002557
$ws.Range("L3").Value = 4.107489981984198
$ws.Range("M3").Value = "[-3.256053817683698, 11.471033781652094]"
$ws.Range("N3").Value = 0.2671851356721482
$ws.Range("O3").Value = 0.2671851356721482
$ws.Range("P3").Value = -0.9182633181663862
$ws.Range("Q3").Value = "[-4.050421759583237, 2.213895123250465]"
$ws.Range("R3").Value = 0.5578227987402717
$ws.Range("S3").Value = 0.5578227987402717
$ws.Range("T3").Value = 10.57708117518096
$ws.Range("U3").Value = "[6.739048372692, 14.415113977669929]"
$ws.Range("V3").Value = [double]"1.443404833878503e-06"
$ws.Range("W3").Value = [double]"1.443404833878503e-06"
$ws.Range("X3").Value = 3.215215215215217
$ws.Range("Y3").Value = -7.751751751751751
$ws.Range("Z3").Value = 14.18218218218218
$ws.Range("H4").Value = 0.2784857805085191
$ws.Range("I4").Value = 0.2784857805085191
$ws.Range("J4").Value = 0.001371436079817601
$ws.Range("K4").Value = 0.001371436079817601
$ws.Range("L4").Value = 3.972531784161972
$ws.Range("M4").Value = "[-2.364706768121327, 10.309770336445272]"
$ws.Range("N4").Value = 0.2132541625149693
$ws.Range("O4").Value = 0.2132541625149693
$ws.Range("P4").Value = -0.1635263443310002
$ws.Range("Q4").Value = "[-3.2579479370560818, 2.9308952483940813]"
$ws.Range("R4").Value = 0.91570923280726
$ws.Range("S4").Value = 0.91570923280726
$ws.Range("T4").Value = 9.215351047281569
$ws.Range("U4").Value = "[5.592913282415344, 12.837788812147794]"
$ws.Range("V4").Value = [double]"6.084682913520822e-06"
$ws.Range("W4").Value = [double]"6.084682913520822e-06"
$ws.Range("X4").Value = 0.5725725725725752
$ws.Range("Y4").Value = -10.26226226226226
$ws.Range("Z4").Value = 11.40740740740741
$ws.Range("F5").Value = 23.22000000000019
$ws.Range("H5").Value = 0.126222436077484
$ws.Range("I5").Value = 0.126222436077484
$ws.Range("L5").Value = 5.741749377629045
$ws.Range("M5").Value = "[-1.1781016217671247, 12.661600377025215]"
$ws.Range("N5").Value = 0.1016228379346555
$ws.Range("O5").Value = 0.1016228379346555
$ws.Range("P5").Value = 0.08176317216550011
$ws.Range("Q5").Value = "[-1.534631846798618, 1.6981581911296182]"
$ws.Range("R5").Value = 0.9193039837615593
$ws.Range("S5").Value = 0.9193039837615593
$ws.Range("T5").Value = 8.669006341723284
$ws.Range("U5").Value = "[4.6850640694187256, 12.652948614027842]"
$ws.Range("V5").Value = [double]"6.95050522845353e-05"
$ws.Range("W5").Value = [double]"6.95050522845353e-05"
$ws.Range("X5").Value = 22.91783783783803
$ws.Range("Y5").Value = 16.94432432432446
$ws.Range("Z5").Value = 28.89135135135159
$ws.Range("F6").Value = 23.22000000000019
$ws.Range("H6").Value = 0.1482887594696771
$ws.Range("I6").Value = 0.1482887594696771
$ws.Range("L6").Value = 6.2055146447282
$ws.Range("M6").Value = "[-2.111143322915198, 14.522172612371598]"
$ws.Range("N6").Value = 0.1398679898862385
$ws.Range("O6").Value = 0.1398679898862385
$ws.Range("P6").Value = -0.2767368904063083
$ws.Range("Q6").Value = "[-2.471763589310889, 1.9182898084982725]"
$ws.Range("R6").Value = 0.8007085769749671
$ws.Range("S6").Value = 0.8007085769749671
$ws.Range("T6").Value = 11.40848133502589
$ws.Range("U6").Value = "[6.918279665161648, 15.898683004890131]"
$ws.Range("V6").Value = [double]"6.217918690376578e-06"
$ws.Range("W6").Value = [double]"6.217918690376578e-06"
$ws.Range("X6").Value = 1.022702702702713
$ws.Range("Y6").Value = -7.089189189189248
$ws.Range("Z6").Value = 9.134594594594674
$ws.Range("F7").Value = 23.22000000000019
$ws.Range("H7").Value = 0.07553616276333963
$ws.Range("I7").Value = 0.07553616276333963
$ws.Range("L7").Value = 6.316036119683565
$ws.Range("M7").Value = "[-0.09922076428176396, 12.731293003648894]"
$ws.Range("N7").Value = 0.05349642164538304
$ws.Range("O7").Value = 0.05349642164538304
$ws.Range("P7").Value = 0.006289474781961069
$ws.Range("Q7").Value = "[-1.3899739268135027, 1.4025528763774249]"
$ws.Range("R7").Value = 0.9928013641075455
$ws.Range("S7").Value = 0.9928013641075455
$ws.Range("T7").Value = 10.60763508931529
$ws.Range("U7").Value = "[6.702635508576119, 14.51263467005447]"
$ws.Range("V7").Value = [double]"1.889201669103002e-06"
$ws.Range("W7").Value = [double]"1.889201669103002e-06"
$ws.Range("X7").Value = 23.19675675675695
$ws.Range("Y7").Value = 18.0367567567569
$ws.Range("Z7").Value = 28.35675675675699
$ws.Range("F8").Value = 23.22000000000019
$ws.Range("H8").Value = 0.01817816738469902
$ws.Range("I8").Value = 0.01817816738469902
$ws.Range("L8").Value = 7.513704881316851
$ws.Range("M8").Value = "[1.382806295544249, 13.644603467089453]"
$ws.Range("N8").Value = 0.01743778683081598
$ws.Range("O8").Value = 0.01743778683081598
$ws.Range("P8").Value = 0.1446579199851161
$ws.Range("Q8").Value = "[-0.9182633181663862, 1.2075791581366184]"
$ws.Range("R8").Value = 0.7852554306129498
$ws.Range("S8").Value = 0.7852554306129498
$ws.Range("T8").Value = 11.05202443957977
$ws.Range("U8").Value = "[7.38706060649328, 14.71698827266627]"
$ws.Range("V8").Value = [double]"2.426787362175276e-07"
$ws.Range("W8").Value = [double]"2.426787362175276e-07"
$ws.Range("X8").Value = 22.68540540540559
$ws.Range("Y8").Value = 18.75729729729744
$ws.Range("Z8").Value = 26.61351351351373
$ws.Range("F9").Value = 23.22000000000019
$ws.Range("H9").Value = 0.08775373781642992
$ws.Range("I9").Value = 0.08775373781642992
$ws.Range("L9").Value = 5.636153714498886
$ws.Range("M9").Value = "[-0.556356543952722, 11.828663972950494]"
$ws.Range("N9").Value = 0.07340054481110614
$ws.Range("O9").Value = 0.07340054481110614
$ws.Range("P9").Value = 0.1949737182408082
$ws.Range("Q9").Value = "[-1.4340002502872329, 1.8239476867688493]"
$ws.Range("R9").Value = 0.8105956835049362
$ws.Range("S9").Value = 0.8105956835049362
$ws.Range("T9").Value = 10.23871344611639
$ws.Range("U9").Value = "[6.650397598388068, 13.82702929384471]"
$ws.Range("V9").Value = [double]"7.404979918135268e-07"
$ws.Range("W9").Value = [double]"7.404979918135268e-07"
$ws.Range("X9").Value = 22.49945945945964
$ws.Range("Y9").Value = 16.47945945945959
$ws.Range("Z9").Value = 28.5194594594597
$ws.Range("F10").Value = 23.22000000000019
$ws.Range("H10").Value = 0.1093218101275631
$ws.Range("I10").Value = 0.1093218101275631
$ws.Range("L10").Value = 6.202342835494605
$ws.Range("M10").Value = "[-1.7101537673743934, 14.114839438363603]"
$ws.Range("N10").Value = 0.1213891039990893
$ws.Range("O10").Value = 0.1213891039990893
$ws.Range("P10").Value = 0.4968685077749626
$ws.Range("Q10").Value = "[-1.6855792415656952, 2.6793162571156204]"
$ws.Range("R10").Value = 0.6487698813158338
$ws.Range("S10").Value = 0.6487698813158338
$ws.Range("T10").Value = 9.409774038990411
$ws.Range("U10").Value = "[5.287333572179609, 13.532214505801214]"
$ws.Range("V10").Value = [double]"3.469243172582459e-05"
$ws.Range("W10").Value = [double]"3.469243172582459e-05"
$ws.Range("X10").Value = 21.38378378378396
$ws.Range("Y10").Value = 13.31837837837848
$ws.Range("Z10").Value = 29.44918918918943
